# "se agregó la capitalizacion." — convert the compounding-growth formulas in
# columns E (cantidad) and F (precio), rows 4-9, of the "ALL" sheet into plain
# capitalizaciones (literal values) with new data, set E3/F3 to their new
# literal values as well, and clear row 10 (E10:F10) entirely so the series
# stops at day 7. All downstream formulas (G, I10, J*, E17, E18, G17, G18, B8…)
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALL")
$ws.Activate()

# Row 3 — literal values already, just new numbers.
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 4.302

# Row 4 — was a formula (=E3*(1+$B$6/100)) / (=F3*(1-$B$5/100)); now literal.
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = 4.229

# Row 5
$ws.Range("E5").Value = 11
$ws.Range("F5").Value = 4.156

# Row 6
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = 4.083

# Row 7
$ws.Range("E7").Value = 19
$ws.Range("F7").Value = 4.009

# Row 8
$ws.Range("E8").Value = 25
$ws.Range("F8").Value = 3.936

# Row 9
$ws.Range("E9").Value = 33
$ws.Range("F9").Value = 3.863

# Row 10 — formulas removed entirely, cells now blank.
$ws.Range("E10").ClearContents()
$ws.Range("F10").ClearContents()

# Selection / scroll moved from C12 to I22 (top-left scrolled to column D).
$ws.Range("I22").Select()

$excel.Calculate()
